$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric need to be pinned to Text format
# first, otherwise Excel auto-converts them to actual numbers on assignment.
$textCells = @('D5','D6','D7','D9','D10','D11','D12','D13','D14','D15','D16','D20','D21','D22','D25','D26','D27','D29','D30','D32','D33','D34','D35','D36','D37','D38','D39','D40','D41','D42','D43','D44','D45','D46','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (prices, % changes, and the TRON/Polkadot row swap)
$ws.Range('D2').Value = '43.927.71'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '2.364.38'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '0.671'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = '239.58'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').Value = '73.93'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').Value = '0.601'
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').Value = '0.102'
$ws.Range('E10').Value = '  +1.84%  '
$ws.Range('D11').Value = '59.77'
$ws.Range('E11').Value = '  +4.32%  '
$ws.Range('D12').Value = '36.83'
$ws.Range('E12').Value = '  +13.69%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.108'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '7.27'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').Value = '16.33'
$ws.Range('E15').Value = '  -1.28%  '
$ws.Range('D16').Value = '0.924'
$ws.Range('E16').Value = '  +2.98%  '
$ws.Range('D17').Value = '2.385.76'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('D18').Value = '43.910.42'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('D20').Value = '6.61'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').Value = '77.49'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').Value = '254.00'
$ws.Range('E22').Value = '  -2.13%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  +3.79%  '
$ws.Range('D25').Value = '1.88'
$ws.Range('E25').Value = '  -3.92%  '
$ws.Range('D26').Value = '2.49'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').Value = '10.57'
$ws.Range('E27').Value = '  -1.34%  '
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('D29').Value = '22.36'
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('D30').Value = '175.54'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('E31').Value = '  +0.97%  '
$ws.Range('D32').Value = '0.134'
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('D33').Value = '0.0758'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').Value = '5.41'
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('D35').Value = '5.08'
$ws.Range('E35').Value = '  -2.56%  '
$ws.Range('D36').Value = '3.79'
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('D37').Value = '6.61'
$ws.Range('E37').Value = '  +5.06%  '
$ws.Range('D38').Value = '2.40'
$ws.Range('E38').Value = '  +2.33%  '
$ws.Range('D39').Value = '0.0282'
$ws.Range('E39').Value = '  +1.57%  '
$ws.Range('D40').Value = '5.54'
$ws.Range('E40').Value = '  +18.21%  '
$ws.Range('D41').Value = '20.51'
$ws.Range('E41').Value = '  +8.39%  '
$ws.Range('D42').Value = '65.95'
$ws.Range('E42').Value = '  +14.24%  '
$ws.Range('D43').Value = '0.107'
$ws.Range('E43').Value = '  -3.06%  '
$ws.Range('D44').Value = '0.203'
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('D45').Value = '9.06'
$ws.Range('E45').Value = '  +1.55%  '
$ws.Range('D46').Value = '2.58'
$ws.Range('E46').Value = '  +2.87%  '
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('D50').Value = '98.26'
$ws.Range('E50').Value = '  -1.61%  '
$ws.Range('D51').Value = '4.39'
$ws.Range('E51').Value = '  +16.23%  '
